$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated measurement values (column G = parent_alkyl, column H = HMW_LMW)
# for rows 2-72, reflecting corrected source-data recalculation ("Sam's updates").
$updates = @{
    "G2" = 2.1
    "H2" = 3.6
    "G3" = 0.2
    "H3" = 1.1000000000000001
    "G4" = 0.6
    "H4" = 1.4
    "G5" = 1.7
    "H5" = 5.7
    "G6" = 2
    "H6" = 6.1
    "G7" = 4.8
    "H7" = 7.1
    "G8" = 3
    "H8" = 7.3
    "G9" = 3.6
    "H9" = 8.5
    "G10" = 1.7
    "H10" = 7.2
    "G11" = 3.1
    "H11" = 10.1
    "G12" = 1.8
    "H12" = 6.9
    "G13" = 2
    "H13" = 5.6
    "G14" = 0.3
    "G15" = 4.0999999999999996
    "H15" = 4.5999999999999996
    "G16" = 0.4
    "H16" = 1.4
    "G17" = 2.5
    "H17" = 5.3
    "G18" = 4.2
    "H18" = 5.8
    "G19" = 2.1
    "H19" = 4.5999999999999996
    "G20" = 3.8
    "H20" = 5.0999999999999996
    "G21" = 3.8
    "H21" = 12.1
    "G22" = 3.7
    "H22" = 6.1
    "G23" = 1.8
    "H23" = 6.1
    "G24" = 2.1
    "H24" = 6
    "G25" = 4.7
    "H25" = 4.5
    "G26" = 0.5
    "H26" = 1.6
    "G27" = 3.4
    "H27" = 4.5999999999999996
    "G28" = 1
    "H28" = 3.6
    "G29" = 3.8
    "H29" = 4.5
    "G30" = 0.7
    "H30" = 3.5
    "G31" = 2.6
    "H31" = 5.0999999999999996
    "G32" = 3.8
    "H32" = 5.4
    "G33" = 1.1000000000000001
    "H33" = 4.3
    "G34" = 2.4
    "H34" = 5.4
    "G35" = 2.5
    "H35" = 10
    "G36" = 4.0999999999999996
    "H36" = 8
    "G37" = 2.4
    "H37" = 5.7
    "G38" = 3.5
    "H38" = 1.1000000000000001
    "G39" = 2.4
    "H39" = 7.4
    "G40" = 2.9
    "H40" = 3.9
    "G41" = 0.2
    "G42" = 1.1000000000000001
    "H42" = 4
    "G44" = 0.6
    "H44" = 0.9
    "G45" = 0.3
    "H45" = 0.6
    "G46" = 0.4
    "H46" = 1.8
    "G47" = 2.1
    "H47" = 4.0999999999999996
    "G48" = 0.3
    "H48" = 0.5
    "G49" = 0.2
    "G50" = 0.9
    "H50" = 1.1000000000000001
    "G51" = 2.9
    "H51" = 13.9
    "B52" = 7
    "G52" = 2.4
    "H52" = 13.9
    "G53" = 2.2999999999999998
    "H53" = 7.6
    "G54" = 2.6
    "H54" = 11.2
    "G55" = 3.5
    "H55" = 5.0999999999999996
    "G56" = 1.8
    "H56" = 4.5
    "G57" = 0.5
    "H57" = 1.4
    "G58" = 2.9
    "H58" = 5.9
    "G59" = 3.5
    "H59" = 5.9
    "G60" = 1.2
    "H60" = 4.2
    "G61" = 3.4
    "H61" = 6.4
    "G62" = 3.4
    "H62" = 5.6
    "G63" = 2
    "H63" = 4.5999999999999996
    "G64" = 2.4
    "H64" = 5.6
    "G65" = 2.7
    "H65" = 8
    "G66" = 3.6
    "H66" = 11.9
    "G67" = 2.8
    "H67" = 4.9000000000000004
    "G68" = 2.2999999999999998
    "H68" = 7
    "G69" = 2.9
    "H69" = 7.3
    "G70" = 1.9
    "H70" = 8
    "G71" = 4
    "H71" = 5.6
    "G72" = 2
    "H72" = 5.7
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value2 = $updates[$key]
}

# B52 unique_id group's sum_EPA16 correction
$ws.Range("B52").Value2 = 7

# Remove the now-obsolete summary/statistics block (rows 73-84: blank separator
# row plus min/median/mean/sd/max/pct>1/pct<1 rollups) so the sheet only holds
# the per-site data table.
$ws.Range("A73:H84").EntireRow.Delete()

# Collapse the leftover split/frozen pane view from the prior session.
$win = $wb.Windows.Item(1)
$win.SplitRow = 0
$win.SplitColumn = 0
$win.Split = $false
$win.FreezePanes = $false

Write-Output $ws.UsedRange.Address()
